$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.777225375175476
$ws.Range("B1").Value = 2.451827526092529
$ws.Range("C1").Value = 2.676681995391846
$ws.Range("D1").Value = 3.463078022003174
$ws.Range("E1").Value = 1.238998293876648
